$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "67÷2=33, 1"
$t.Cell(1,2).Range.Text = "68÷7=9, 5"
$t.Cell(1,3).Range.Text = "88÷6=14, 4"
$t.Cell(1,4).Range.Text = "48÷9=5, 3"
$t.Cell(1,5).Range.Text = "76÷8=9, 4"
$t.Cell(5,1).Range.Text = "76÷6=12, 4"
$t.Cell(5,2).Range.Text = "93÷6=15, 3"
$t.Cell(5,3).Range.Text = "31÷2=15, 1"
$t.Cell(5,4).Range.Text = "63÷5=12, 3"
$t.Cell(5,5).Range.Text = "47÷8=5, 7"
$t.Cell(9,1).Range.Text = "67÷7=9, 4"
$t.Cell(9,2).Range.Text = "78÷9=8, 6"
$t.Cell(9,3).Range.Text = "79÷6=13, 1"
$t.Cell(9,4).Range.Text = "84÷8=10, 4"
$t.Cell(9,5).Range.Text = "66÷4=16, 2"
$t.Cell(13,1).Range.Text = "77÷7=11, 0"
$t.Cell(13,2).Range.Text = "31÷6=5, 1"
$t.Cell(13,3).Range.Text = "30÷9=3, 3"
$t.Cell(13,4).Range.Text = "19÷4=4, 3"
$t.Cell(13,5).Range.Text = "98÷5=19, 3"
$t.Cell(17,1).Range.Text = "69÷6=11, 3"
$t.Cell(17,2).Range.Text = "99÷8=12, 3"
$t.Cell(17,3).Range.Text = "85÷8=10, 5"
$t.Cell(17,4).Range.Text = "63÷2=31, 1"
$t.Cell(17,5).Range.Text = "98÷7=14, 0"
